# Add "Driver" column (J) to the Trip report template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the jx: template comments so their lastCell points at the new
#    last column (J9 instead of I9).
$ws.Comments.Item(1).Text('jx:area(lastCell="J9")')
$ws.Comments.Item(2).Text("jx:each(items=`"devices`", var=`"device`", lastCell=`"J9`" multisheet=`"sheetNames`")`n")
$ws.Comments.Item(3).Text('jx:each(items="device.objects", var="trip", lastCell="J9")')

# 2. Add the new header cell (J8) and template expression cell (J9),
#    copying formatting from the existing last column (I) so the new
#    column matches the look of the table.
$ws.Range("I8").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("J8").Value = "Driver"

$ws.Range("I9").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("J9").Value = '${trip.driverName ? trip.driverName : trip.driverUniqueId}'

# 3. Give column J its own width (matches the widened column in the
#    authored template; ColumnWidth is quantised to the nearest 1/6
#    character by this engine, 16.667 is the closest achievable value).
$ws.Columns.Item(10).ColumnWidth = 16.66666666666667
